# English presentation STARRT - text revision pass
# Applies the wording/content fixes from the commit (task #79).

$d = $word.ActiveDocument

# Common quote characters used throughout this document.
$lsq = [char]0x2018   # '
$rsq = [char]0x2019   # '
$lbrk = [char]11      # manual line break (w:br) as it appears in Range.Text

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Output "NOT FOUND: $find"
    }
}

# 1. "Goodmorning" -> "Good morning"
Replace-Text "Goodmorning everybody, as " "Good morning everybody, as "

# 2. "As IT-students" -> "As IT students"
Replace-Text "As IT-students of PXL, w" "As IT students of PXL, w"

# 3. remove "cross-department "
Replace-Text "This platform will be used cross-department to provide" "This platform will be used to provide"

# 4. "OfficeCenter or various locations." -> "OfficeCenter or various other locations."
Replace-Text "OfficeCenter or various locations. " "OfficeCenter or various other locations. "

# 5. "a pure ideas-pitcher." -> "purely an ideas pitcher."
Replace-Text "get involved. The platform is a pure ideas-pitcher." "get involved. The platform is purely an ideas pitcher."

# 6. "Dennie was mostly into security" -> "Dennie set up a basic Spring Boot backend and some security"
Replace-Text "Dennie was mostly into security" "Dennie set up a basic Spring Boot backend and some security"

# 7. Big paragraph: "opensource" -> "open source", "creating a plan" -> "forming a plan",
#    "more fuzzy than we first tought" -> "fuzzier than we first thought",
#    "So this plan is in current development and we dediced to go" -> "This initial plan is currently still in development, but we decided to go"
$find7 = "Our goal was to create a web-based future-proof and secure Single Page Application with opensource libraries and new techniques. We started the creation of the platform by creating a plan. This plan was obviously intended to meet the requirements but we soon realized that the requirements were more fuzzy than we first tought. So this plan is in current development and we dediced to go " + $lsq + "Agile" + $rsq + "."
$repl7 = "Our goal was to create a web-based future-proof and secure Single Page Application with open source libraries and new techniques. We started the creation of the platform by forming a plan. This plan was obviously intended to meet the requirements but we soon realized that the requirements were fuzzier than we first thought. This initial plan is currently still in development, but we decided to go " + $lsq + "Agile" + $rsq + "."
Replace-Text $find7 $repl7

# 8. brainstorm paragraph
$find8 = "As a first step, we brainstormed a lot. We created an initial analysis and were surprised to see that everyone had a different idea on how we should move on. So we decided to arrange some more meeting to discuss how we could meet the goals that were requested. "
$repl8 = "As a first step, we brainstormed a lot. We created an initial exploring analysis and realized that everyone had a different idea on how we should move on. So we decided to arrange some more meeting to discuss how we could meet the goals that were requested and at the same time combine our visions. "
Replace-Text $find8 $repl8

# 9. "Next, we started the analysis process" -> "Next, we started the deeper analysis process"
Replace-Text "Next, we started the analysis process" "Next, we started the deeper analysis process"

# 10. "for example: how shall we upload" -> "For example: how do we upload"
Replace-Text "for example: how shall we upload" "For example: how do we upload"

# 11. "using skype, everyday at" -> "using Skype, every day at"
Replace-Text "using skype, everyday at" "using Skype, every day at"

# 12. "have a job it's not easy to get into all sessions" -> "have a full time job it's not easy to attend all sessions"
$find12 = "Since most of us have a job it" + $rsq + "s not easy to get into all sessions"
$repl12 = "Since most of us have a full time job it" + $rsq + "s not easy to attend all sessions"
Replace-Text $find12 $repl12

# 13. "When analysis was mostly clear, we started the initial development. Everyone using their own skills and expertise we managed to get a potentially shippable product"
$find13 = "When analysis was mostly clear, we started the initial development. Everyone using their own skills and expertise we managed to get a potentially shippable product"
$repl13 = "When the biggest part of our analysis was done, we started the initial development. Using everyone" + $rsq + "s own skills and expertise we managed to get a potentially shippable product"
Replace-Text $find13 $repl13

# 14. "If course, this product" -> "Of course, this product"
Replace-Text "time.If course, this product had only about 20pct" "time.Of course, this product had only about 20pct"

# 15. "stories and requirements clear." -> "stories and requirements done."
Replace-Text "we had some user stories and requirements clear." "we had some user stories and requirements done."

# 17. "provided as-is. Accomplishing" -> "provided as-is and only needs fine-tuning. Accomplishing"
Replace-Text "provided as-is. Accomplishing this task was more difficult" "provided as-is and only needs fine-tuning. Accomplishing this task was more difficult"

# 18. "A guest for example," -> "A guest, for example,"
Replace-Text "A guest for example, can only see" "A guest, for example, can only see"

# 19. "this small user stories contain a certain grade of complexity and should be discussed in group." ->
#     "this small user story contains a certain grade of complexity and should be discussed as a group."
Replace-Text "this small user stories contain a certain grade of complexity and should be discussed in group." "this small user story contains a certain grade of complexity and should be discussed as a group."

# 20. " had the possibility to log-in." -> " has the possibility to log in."
Replace-Text "phase 2, when a user had the possibility to log-in." "phase 2, when a user has the possibility to log in."

# 22. "we have to overcome before other tasks could get completed." -> "we had to overcome before other tasks could get completed."
Replace-Text "unplanned obstacles we have to overcome before other tasks could get completed." "unplanned obstacles we had to overcome before other tasks could get completed."

# 23. Meetings paragraph full rewrite.
$find23 = "We held lots of online-meetings. We tend to discuss lots of " + $lsq + "irrelevant" + $rsq + " information during these meetings since we all want to go into detail. These should get more to-the-point since development-time gets lost during long meetings where just a little is descided."
$repl23 = "We held lots of online-meetings. We tend to discuss lots of " + $lsq + "irrelevant" + $rsq + " information during these meetings since we all wanted to go into detail. These should have been more to the point since development-time gets lost during long meetings where too little is decided."
Replace-Text $find23 $repl23

# 25. "Since face-to-face..." paragraph: pluralize "meetings", drop the manual line break
#     (replace with a space), "since these are" -> "since they are", "Phone" -> "phone".
$find25 = "Since face-to-face is the most effective communication, we decided to plan some meeting at OffiCenter." + $lbrk + "During these meetings, most of the requirements were discussed since these are the foundations of our software. During Phone conversations we held every day, we discussed the progress & current issues."
$repl25 = "Since face-to-face is the most effective communication, we decided to plan some meetings at OffiCenter. During these meetings, most of the requirements were discussed since they are the foundations of our software. During phone conversations we held every day, we discussed the progress & current issues."
Replace-Text $find25 $repl25

# 26. "if course,d evelopment" -> "of course, development"; "Webstorm while Jesse" -> "Webstorm, while Jesse"
Replace-Text "Because, if course,d evelopment could go much faster when using tools we know." "Because, of course, development could go much faster when using tools we know."
Replace-Text "For example, Arjen likes to work with Webstorm while Jesse prefers" "For example, Arjen likes to work with Webstorm, while Jesse prefers"

# 27. Append " with a broad community" after "... open source technology" and move the
#     _GoBack bookmark there (Word always keeps _GoBack at the location of the most
#     recent text edit).
Replace-Text "Known & supported open source technology" "Known & supported open source technology with a broad community"

try {
    $oldBookmark = $d.Bookmarks("_GoBack")
    $oldBookmark.Delete()
} catch {
}

$found = $d.Content.Find.Execute("technology with a broad community", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $endRange = $d.Range($d.Content.Find.Parent.End, $d.Content.Find.Parent.End)
    $d.Bookmarks.Add("_GoBack", $endRange)
}

Write-Output "Done"
